$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Work bottom-up (highest original row number first) so that row numbers
# referenced below (taken from the *original* layout before any inserts)
# stay valid as each insertion is performed.

# 4) Insert a new row right after original row 7 (LEONARDO, 2609.8) -> new row 8: DIOGO
$ws.Rows.Item(8).Insert()
$ws.Range("A8").NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "004550415"
$ws.Cells.Item(8,2).Value = "DIOGO"
$ws.Cells.Item(8,3).Value = 2155.74

# 3) Insert 5 new rows right before original row 5 (GUSTAVO) -> rows 5-9
$ws.Range("A5:A9").EntireRow.Insert()

$ws.Range("A5").NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "005061124"
$ws.Cells.Item(5,2).Value = "BRUNO"
$ws.Cells.Item(5,3).Value = 14309.72

$ws.Range("A6").NumberFormat = "@"
$ws.Cells.Item(6,1).Value = "005547702"
$ws.Cells.Item(6,2).Value = "NATHALIA"
$ws.Cells.Item(6,3).Value = 9960.05

$ws.Range("A7").NumberFormat = "@"
$ws.Cells.Item(7,1).Value = "004515548"
$ws.Cells.Item(7,2).Value = "FLAVIA"
$ws.Cells.Item(7,3).Value = 9907.13

$ws.Range("A8").NumberFormat = "@"
$ws.Cells.Item(8,1).Value = "004886366"
$ws.Cells.Item(8,2).Value = "RENATO"
$ws.Cells.Item(8,3).Value = 8806.33

$ws.Range("A9").NumberFormat = "@"
$ws.Cells.Item(9,1).Value = "004202332"
$ws.Cells.Item(9,2).Value = "TATIANA"
$ws.Cells.Item(9,3).Value = 6604.48

# 2) Insert 1 new row right before original row 4 (KAUANNE) -> new row 4: CLISIA
$ws.Rows.Item(4).Insert()
$ws.Range("A4").NumberFormat = "@"
$ws.Cells.Item(4,1).Value = "004805273"
$ws.Cells.Item(4,2).Value = "CLISIA"
$ws.Cells.Item(4,3).Value = 23115.69

# 1) Overwrite row 2 in place (INTERLAGOS/272337.4 -> BRUNO/51735.12)
$ws.Range("A2").NumberFormat = "@"
$ws.Cells.Item(2,1).Value = "004515341"
$ws.Cells.Item(2,2).Value = "BRUNO"
$ws.Cells.Item(2,3).Value = 51735.12
